$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New component names (interned into shared strings first, in row order)
$ws.Range("A6").Value = "Buck Boost"
$ws.Range("A7").Value = "Boost"

# New links. The "Boost" row's link (TPS55340RTER) is interned before the
# "Buck Boost" row's link (TPS630702RNMR) to match the shared-string order.
$ws.Range("B7").Value = "https://www.digikey.ca/en/products/detail/texas-instruments/TPS55340RTER/3503781"
$ws.Range("B6").Value = "https://www.digikey.ca/en/products/detail/texas-instruments/TPS630702RNMR/10434765"

# Mech/Elec/Firm classification column (reuses existing "Elec" shared string)
$ws.Range("C6").Value = "Elec"
$ws.Range("C7").Value = "Elec"

# Update the selected cell to match the saved workbook state
$ws.Range("E5").Select()
